# kidum.xlsx - "Add files via upload" edit
#
# What changed (per the OOXML diff):
#   1. The shared string used by A6 ("free games helix jump") had stray
#      leading/trailing whitespace trimmed off.
#   2. A4:B6 (which had a "plain" look) were reformatted to match the same
#      font/wrap formatting already used by the header block A1:B3 - this
#      also happens to retire the now-unused "plain + wrap" cell format.
#   3. Row 6 no longer needs its tall custom row height (it used to be
#      taller to fit wrapped text in the old font) - it goes back to the
#      sheet's default row height.
#   4. Column A was widened.
#   5. A new, empty, similarly formatted row 7 was added below the table,
#      and the sheet's selection now sits on A7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clean up the "free games helix jump" keyword text -----------------
$ws.Range("A6").Value = "free games helix jump"

# --- 2. Re-apply the header row's formatting down over A4:B7 --------------
# (copies the cell format only - values/formulas are left alone - and
# reuses the existing font/wrap format instead of inventing a new one)
$ws.Range("A1:B1").Copy()
$ws.Range("A4:B7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# --- 3. Row 6 returns to the sheet's default (non-custom) row height ------
$ws.Rows.Item(6).AutoFit()

# --- 4. Widen column A ------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 27.2857142857143

# --- 5. Leave the selection on the new blank row ---------------------------
[void]$ws.Range("A7").Select()
